$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 243; this shifts existing rows 243-287 down to 244-288
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across all records in this sheet.
$ws.Cells.Item(243, 1).Value = 3
$ws.Cells.Item(243, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(243, 3).Value = "Coquimbo"
$ws.Cells.Item(243, 4).Value = 44694
$ws.Cells.Item(243, 5).Value = 5
$ws.Cells.Item(243, 6).Value = 100112001
$ws.Cells.Item(243, 7).Value = "Berenjena"
$ws.Cells.Item(243, 8).Value = "Sin especificar"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 115
$ws.Cells.Item(243, 11).Value = 8500
$ws.Cells.Item(243, 12).Value = 9000
$ws.Cells.Item(243, 13).Value = 8739
$ws.Cells.Item(243, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(243, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(243, 16).Value = 146
$ws.Cells.Item(243, 17).Value = 60
$ws.Cells.Item(243, 18).Value = "Hortaliza"
